$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A1").Value = "First Name"
$ws1.Range("B1").Value = "Last Name"
$ws1.Range("D1").Value = "Upload Document"
$ws1.Range("C1").Value = "Comments"

$ws1.Columns("A:D").AutoFit() | Out-Null

$ws1.Range("C1").Select() | Out-Null
